# Natmi following Dr Hou advice
# Rewrites the Il16-Cd4 LR-pair sheet data (rows 2-13) to add the
# previously-missing "self" sending/target cluster combinations, turning
# the 3-target-per-sender grid into a full 4-target-per-sender grid
# (ECs/FAPs/M2/sCs x ECs/FAPs/M2/sCs), rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data block first (A2:T13) so leftover cells from the
# previous, shorter table don't linger once the new, longer table is
# written below.
$ws.Range("A2:T13").ClearContents()

$rows = @(
  @{ A="ECs"; B="Il16"; C="Cd4"; D="ECs"; Nums=@(3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,2,0.6666666666666666,1.033717666666667,3.101153,0.03978942744299105,0.03978942744299105,4.386708410345556,39.48037569311001,0.009824990261906136,0.009824990261906138) },
  @{ A="ECs"; B="Il16"; C="Cd4"; D="FAPs"; Nums=@(3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,3,1,7.139908666666667,21.419726,0.2748263737796068,0.2748263737796068,30.29908301573556,272.6917471416201,0.06786140489124454,0.06786140489124455) },
  @{ A="ECs"; B="Il16"; C="Cd4"; D="M2"; Nums=@(3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,3,1,17.222987,51.668961,0.6629400015943212,0.6629400015943213,73.08786950289667,657.7908255260701,0.1636962248130963,0.1636962248130963) },
  @{ A="ECs"; B="Il16"; C="Cd4"; D="sCs"; Nums=@(3,1,4.243623333333334,12.73087,0.2469246453968972,0.2469246453968973,3,1,0.5830936666666667,1.749281,0.02244419718308088,0.02244419718308088,2.474429889385556,22.26986900447,0.005542025430650287,0.005542025430650288) },
  @{ A="FAPs"; B="Il16"; C="Cd4"; D="ECs"; Nums=@(3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,2,0.6666666666666666,1.033717666666667,3.101153,0.03978942744299105,0.03978942744299105,6.063262359519445,54.569361235675,0.01357999850119236,0.01357999850119236) },
  @{ A="FAPs"; B="Il16"; C="Cd4"; D="FAPs"; Nums=@(3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,3,1,7.139908666666667,21.419726,0.2748263737796068,0.2748263737796068,41.87907478509445,376.91167306585,0.09379732214952019,0.09379732214952023) },
  @{ A="FAPs"; B="Il16"; C="Cd4"; D="M2"; Nums=@(3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,3,1,17.222987,51.668961,0.6629400015943212,0.6629400015943213,101.0212867236083,909.1915805124752,0.2262592051853509,0.226259205185351) },
  @{ A="FAPs"; B="Il16"; C="Cd4"; D="sCs"; Nums=@(3,1,5.865491666666667,17.596475,0.3412966552647515,0.3412966552647516,3,1,0.5830936666666667,1.749281,0.02244419718308088,0.02244419718308088,3.420131042719445,30.781179384475,0.007660129428688063,0.007660129428688065) },
  @{ A="M2"; B="Il16"; C="Cd4"; D="ECs"; Nums=@(3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,2,0.6666666666666666,1.033717666666667,3.101153,0.03978942744299105,0.03978942744299105,6.329569049378666,56.96612144440799,0.01417645041679648,0.01417645041679648) },
  @{ A="M2"; B="Il16"; C="Cd4"; D="FAPs"; Nums=@(3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,3,1,7.139908666666667,21.419726,0.2748263737796068,0.2748263737796068,43.71846043577066,393.466143921936,0.09791702749924508,0.0979170274992451) },
  @{ A="M2"; B="Il16"; C="Cd4"; D="M2"; Nums=@(3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,3,1,17.222987,51.668961,0.6629400015943212,0.6629400015943213,105.458278375544,949.1245053798959,0.2361968157339837,0.2361968157339838) },
  @{ A="M2"; B="Il16"; C="Cd4"; D="sCs"; Nums=@(3,1,6.123111999999999,18.369336,0.3562868663317164,0.3562868663317164,3,1,0.5830936666666667,1.749281,0.02244419718308088,0.02244419718308088,3.570347827490667,32.13313044741599,0.007996572681691023,0.007996572681691025) },
  @{ A="sCs"; B="Il16"; C="Cd4"; D="ECs"; Nums=@(3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,2,0.6666666666666666,1.033717666666667,3.101153,0.03978942744299105,0.03978942744299105,0.9858331077662223,8.872497969896001,0.002207988263096068,0.002207988263096068) },
  @{ A="sCs"; B="Il16"; C="Cd4"; D="FAPs"; Nums=@(3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,3,1,7.139908666666667,21.419726,0.2748263737796068,0.2748263737796068,6.809169057470223,61.28252151723201,0.01525061923959691,0.01525061923959692) },
  @{ A="sCs"; B="Il16"; C="Cd4"; D="M2"; Nums=@(3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,3,1,17.222987,51.668961,0.6629400015943212,0.6629400015943213,16.42517231419467,147.826550827752,0.03678775586189022,0.03678775586189024) },
  @{ A="sCs"; B="Il16"; C="Cd4"; D="sCs"; Nums=@(3,1,0.9536773333333334,2.861032,0.05549183300663471,0.05549183300663472,3,1,0.5830936666666667,1.749281,0.02244419718308088,0.02244419718308088,0.5560832131102224,5.004748917992001,0.001245469642051506,0.001245469642051506) }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    $c = 5
    foreach ($n in $row.Nums) {
        $ws.Cells.Item($r, $c).Value = $n
        $c++
    }

    $r++
}
